# Weatherdata.xlsx edit: add a "Condition" column (H) derived from a
# tree-based weather classifier, overwrite the C/E/F/G sample values with
# the model's feature values, and trim the sheet down to rows 2-23
# (removing the trailing duplicate rows 24-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data for rows 2-23 ------------------------------------------------
# Columns: Day(C), Month(D), Temprature(E), Humidity(F), WindSpeed(G), Condition(H)
$data = @(
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 30, 37,  6.17, "['Hot', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 29, 37,  6.17, "['Hot', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 100, 6.17, "['Cold', 'High Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 75,  6.17, "['Cold', 'High Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 10, 10,  10,   "['Cold', 'Low Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  6.17, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']"),
    @(2, 3, 21, 37,  5.14, "['Cold', 'Mild Humidity', 'Low Wind']")
)

# --- Header for the new Condition column -----------------------------------
$ws.Cells.Item(1, 8).Value = "Condition"

# --- Write the updated rows 2-23 -------------------------------------------
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $r++
}

# --- Remove the now-unused trailing rows 24-28 ------------------------------
$ws.Range("A24:H28").EntireRow.Delete()

# --- Column widths (best-fit, as Excel would auto-size after entry) --------
# Target OOXML widths are 10.21875 / 10.109375 / 59.44140625 characters;
# Excel's ColumnWidth property is offset from the stored sheet width by the
# default font padding (~0.8333 chars for Calibri 11) and then the stored
# value itself is quantized to pixel boundaries on save, so the values below
# are chosen to land on the closest achievable stored width.
$ws.Columns.Item(1).ColumnWidth = 9.333333333333332
$ws.Columns.Item(7).ColumnWidth = 9.333333333333332
$ws.Columns.Item(8).ColumnWidth = 58.666666666666664

# --- View / selection state ---------------------------------------------------
$ws.Range("A7").Select()
$ws.Range("H14").Select()
